$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GLOBAL RESULTS")

$ws.Range("C2").Value = 5.443642818356839
$ws.Range("C3").Value = 12.717880509399059
$ws.Range("C4").Value = 4.271532611314539
$ws.Range("C5").Value = 10.019195499471632
$ws.Range("C6").Value = 4.670094009256409
$ws.Range("C7").Value = 10.93684954232447
$ws.Range("C8").Value = 4.164794848457012
$ws.Range("C9").Value = 9.773440792374878
